$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Biz Widget"
$ws.Range("B5").Value = 400.0
